$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change, label stays the same
$ws.Range("B3").Value = 0.9955494309793238
$ws.Range("C3").Value = 0.9962973768916696
$ws.Range("D3").Value = 0.9960579398049054

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9962983603031225
$ws.Range("C4").Value = 0.996427947679896
$ws.Range("D4").Value = 0.9964287679898766

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9949272161718722
$ws.Range("C5").Value = 0.9951390687363776
$ws.Range("D5").Value = 0.995171874322312
